$d = $word.ActiveDocument

# Update the placeholder ID text in the first paragraph (also absorbs the
# trailing space run that used to follow it, so the paragraph ends up with
# a single run and no extra trailing space).
$d.Content.Find.Execute(
    "**ID__AFFARS_pgi_5304_topic_6__ID** ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "**ID__AFFARS_AFMC_PGI_5304_1602__ID**",
    2
)

# Give the first paragraph the same paragraph border / indent treatment
# used by the rest of the list paragraphs in the document.
$p = $d.Paragraphs(1)
$p.Format.LeftIndent = 11.25

$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

Write-Host "Done"
